$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A (Mat / NC) - numeric ids
$ws.Cells.Item(2, 1).Value = 19330051920159
$ws.Cells.Item(3, 1).Value = 19330051920165

# Column B (Paterno)
$ws.Cells.Item(2, 2).Value = "GARCIA"
$ws.Cells.Item(3, 2).Value = "MACARIO"

# Column C (Materno)
$ws.Cells.Item(2, 3).Value = "ZUÑIGA"
$ws.Cells.Item(3, 3).Value = "NIEVES"

# Column D (Nombres)
$ws.Cells.Item(2, 4).Value = "MARCO ANTONIO"
$ws.Cells.Item(3, 4).Value = "MARTHA"

# Column E (Nombre_Largo)
$ws.Cells.Item(2, 5).Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"
$ws.Cells.Item(3, 5).Value = "ANALIZA Y FRACCIONA SANGRE CON FINES TRANSFUSIONALES"

# Column F (Grupo)
$ws.Cells.Item(2, 6).Value = "5ALCM"
$ws.Cells.Item(3, 6).Value = "5ALCM"

# Column G (Reprobadas)
$ws.Cells.Item(2, 7).Value = 6
$ws.Cells.Item(3, 7).Value = 6
